$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Diffraction-angle helper table (columns P-S), rows 5-10
$ws.Range("P5").Formula = "=E6+1"
$ws.Range("Q5").Formula = "=P5/2"
$ws.Range("R5").Formula = "=SQRT((49.4 - F6)/100 * Q5 * 546.1/1000000000)"
$ws.Range("S5").Formula = "=2*R5"
$ws.Range("P6").Formula = "=E7+1"
$ws.Range("Q6").Formula = "=P6/2"
$ws.Range("R6").Formula = "=SQRT((49.4 - F7)/100 * Q6 * 546.1/1000000000)"
$ws.Range("S6").Formula = "=2*R6"
$ws.Range("P7").Formula = "=E8+1"
$ws.Range("Q7").Formula = "=P7/2"
$ws.Range("R7").Formula = "=SQRT((49.4 - F8)/100 * Q7 * 546.1/1000000000)"
$ws.Range("S7").Formula = "=2*R7"
$ws.Range("P8").Formula = "=E9+1"
$ws.Range("Q8").Formula = "=P8/2"
$ws.Range("R8").Formula = "=SQRT((49.4 - F9)/100 * Q8 * 546.1/1000000000)"
$ws.Range("S8").Formula = "=2*R8"
$ws.Range("P9").Formula = "=E10+1"
$ws.Range("Q9").Formula = "=P9/2"
$ws.Range("R9").Formula = "=SQRT((49.4 - F10)/100 * Q9 * 546.1/1000000000)"
$ws.Range("S9").Formula = "=2*R9"
$ws.Range("P10").Formula = "=E11+1"
$ws.Range("Q10").Formula = "=P10/2"
$ws.Range("R10").Formula = "=SQRT((49.4 - F11)/100 * Q10 * 546.1/1000000000)"
$ws.Range("S10").Formula = "=2*R10"

# Diffraction-angle helper table (columns H-I), rows 6-11
$ws.Range("H6").Formula = "=SQRT((49.4 - F6)/100 * E6 * 5461/10000000000)"
$ws.Range("I6").Formula = "=2*H6"
$ws.Range("H7").Formula = "=SQRT((49.4 - F7)/100 * E7 * 5461/10000000000)"
$ws.Range("I7").Formula = "=2*H7"
$ws.Range("H8").Formula = "=SQRT((49.4 - F8)/100 * E8 * 5461/10000000000)"
$ws.Range("I8").Formula = "=2*H8"
$ws.Range("H9").Formula = "=SQRT((49.4 - F9)/100 * E9 * 5461/10000000000)"
$ws.Range("I9").Formula = "=2*H9"
$ws.Range("H10").Formula = "=SQRT((49.4 - F10)/100 * E10 * 5461/10000000000)"
$ws.Range("I10").Formula = "=2*H10"
$ws.Range("H11").Formula = "=SQRT((49.4 - F11)/100 * E11 * 5461/10000000000)"
$ws.Range("I11").Formula = "=2*H11"

# K/L helper table rows 12-15
$ws.Range("K12").Value = -2
$ws.Range("L12").Value = 28
$ws.Range("K13").Value = -1
$ws.Range("L13").Value = 32
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 102
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 106

# I14, I15:I21 (1e6 * I column), plus standalone C15/C20
$ws.Range("I14").Formula = "=1000000*I6"
$ws.Range("C15").Formula = "=(3.8 - 1.42) * 400"
$ws.Range("I15").Formula = "=1000000*I7"
$ws.Range("I16").Formula = "=1000000*I8"
$ws.Range("I17").Formula = "=1000000*I9"
$ws.Range("I18").Formula = "=1000000*I10"
$ws.Range("I19").Formula = "=1000000*I11"
$ws.Range("I20").Formula = "=1000000*I12"
$ws.Range("I21").Formula = "=1000000*I13"
$ws.Range("C20").Formula = "=1.4 - 2.4"

# Slit-width / diffraction section rows 24, 26
$ws.Range("F24").Value = "D"
$ws.Range("G24").Formula = "=212/1000000"
$ws.Range("I24").Formula = "=12.5/100*546.1/1000000000/G24"

$ws.Range("H26").Value = -1
$ws.Range("I26").Value = -0.3
$ws.Range("J26").Value = 0.4
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 1.7
$ws.Range("M26").Value = 2.5
$ws.Range("N26").Value = 3
$ws.Range("O26").Value = 3.8

# C/D/E/F table rows 28-34 (plus E35)
$ws.Range("D28").Value = -1
$ws.Range("C28").Formula = "=1.5 + D28"
$ws.Range("E28").Formula = "=400*D28"
$ws.Range("F28").Formula = "=400*C28"
$ws.Range("D29").Value = -0.2
$ws.Range("C29").Formula = "=1.5 + D29"
$ws.Range("E29").Formula = "=400*D29"
$ws.Range("F29").Formula = "=400*C29"
$ws.Range("D30").Value = 0.7
$ws.Range("C30").Formula = "=1.5 + D30"
$ws.Range("E30").Formula = "=400*D30"
$ws.Range("F30").Formula = "=400*C30"
$ws.Range("D31").Value = 1.4
$ws.Range("C31").Formula = "=1.5 + D31"
$ws.Range("E31").Formula = "=400*D31"
$ws.Range("F31").Formula = "=400*C31"
$ws.Range("D32").Value = 2.3
$ws.Range("C32").Formula = "=1.5 + D32"
$ws.Range("E32").Formula = "=400*D32"
$ws.Range("F32").Formula = "=400*C32"
$ws.Range("D33").Value = 3.1
$ws.Range("C33").Formula = "=1.5 + D33"
$ws.Range("E33").Formula = "=400*D33"
$ws.Range("F33").Formula = "=400*C33"
$ws.Range("D34").Value = 3.8
$ws.Range("C34").Formula = "=1.5 + D34"
$ws.Range("E34").Formula = "=400*D34"
$ws.Range("F34").Formula = "=400*C34"
$ws.Range("E35").Formula = "=400*D35"

# Column widths for H and R (best-fit-like)
$ws.Columns("H").ColumnWidth = 10.92
$ws.Columns("R").ColumnWidth = 10.92

# Selection matches where the author left off
[void]$ws.Range("S14").Select()
